$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55 (2024-02-23) precipitation value was recorded after the fact: 0 -> 0.8
$ws.Range("B55").Value = 0.8

# Append new daily rows (2024-02-24 .. 2024-05-17) with their precipitation readings
$dates = @("2024-02-24","2024-02-25","2024-02-26","2024-02-27","2024-02-28","2024-02-29","2024-03-01","2024-03-02","2024-03-03","2024-03-04","2024-03-05","2024-03-06","2024-03-07","2024-03-08","2024-03-09","2024-03-10","2024-03-11","2024-03-12","2024-03-13","2024-03-14","2024-03-15","2024-03-16","2024-03-17","2024-03-18","2024-03-19","2024-03-20","2024-03-21","2024-03-22","2024-03-23","2024-03-24","2024-03-25","2024-03-26","2024-03-27","2024-03-28","2024-03-29","2024-03-30","2024-03-31","2024-04-01","2024-04-02","2024-04-03","2024-04-04","2024-04-05","2024-04-06","2024-04-07","2024-04-08","2024-04-09","2024-04-10","2024-04-11","2024-04-12","2024-04-13","2024-04-14","2024-04-15","2024-04-16","2024-04-17","2024-04-18","2024-04-19","2024-04-20","2024-04-21","2024-04-22","2024-04-23","2024-04-24","2024-04-25","2024-04-26","2024-04-27","2024-04-28","2024-04-29","2024-04-30","2024-05-01","2024-05-02","2024-05-03","2024-05-04","2024-05-05","2024-05-06","2024-05-07","2024-05-08","2024-05-09","2024-05-10","2024-05-11","2024-05-12","2024-05-13","2024-05-14","2024-05-15","2024-05-16","2024-05-17")
$values = @(14.2,14.2,3.4,10,19,14,4.6,1,4.4,5.2,8.2,4,2.2,2.4,3.2,7.6,0.6,0.6,5.4,0.6,0,0.2,4.2,0.4,1.2,0,0,0.4,0,0,1.2,0,0.6,2,0,0,0,0,0,0,2.6,0,0.4,4.6,0,0,0,0,0,0,0,0,0,0,0,0,0,1.2,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$startRow = 56
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $dates[$i]
    $cellA.Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $values[$i]
}
